$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in "Points for grading" (column E) scores that graders had left blank ---

# Rubric section 1 (Customer Class), rows 3-6
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("E6").Value = 2
# Section total: hard-typed as 7 rather than re-deriving the SUM formula
$ws.Range("E7").Formula = "=7"

# Rubric section 2 (Product Class), rows 10-14
$ws.Range("E10").Value = 2
$ws.Range("E11").Value = 2
$ws.Range("E12").Value = 2
$ws.Range("E13").Value = 2
$ws.Range("E14").Value = 2
# Row 15's section total (=SUM(E10:E14)) recalculates automatically to 10

# --- Update the active selection / scroll position left by the grader ---
$null = $ws.Range("E15").Select()
